# Fix the Start Date values in the DeliveryChanges sheet so they show as
# plain dates (no time component) and correct the typo'd day in A2.
#
# These cells already contain plain text (not real Excel dates), so before
# typing the corrected values back in we mark the range as Text-formatted.
# That mirrors what a user does via the GUI (Format Cells -> Text) to stop
# Excel from reinterpreting a date-looking string as a serial date number,
# keeping the result a literal string like "2025-04-07".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateRange = $ws.Range("A2:A6")
$dateRange.NumberFormat = "@"

$ws.Range("A2").Value = "2025-04-07"
$ws.Range("A3").Value = "2025-04-07"
$ws.Range("A4").Value = "2025-04-07"
$ws.Range("A5").Value = "2025-04-08"
$ws.Range("A6").Value = "2025-04-09"
